$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute("2023-04-23 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-04-24 Monday", 2) | Out-Null

# Update each multiplication-problem cell by its table position
# (value-based Find/Replace is unsafe because several cells share
# the same "old" text, e.g. "38×40=" and "58×77=" each appear twice
# and must map to different new values).
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "18×36="
$t.Cell(1,2).Range.Text = "90×91="
$t.Cell(1,3).Range.Text = "76×18="
$t.Cell(1,4).Range.Text = "72×82="
$t.Cell(1,5).Range.Text = "46×76="
$t.Cell(2,1).Range.Text = "19×97="
$t.Cell(2,2).Range.Text = "62×94="
$t.Cell(2,3).Range.Text = "86×15="
$t.Cell(2,4).Range.Text = "51×96="
$t.Cell(2,5).Range.Text = "44×70="
$t.Cell(3,1).Range.Text = "61×62="
$t.Cell(3,2).Range.Text = "76×56="
$t.Cell(3,3).Range.Text = "50×92="
$t.Cell(3,4).Range.Text = "12×58="
$t.Cell(3,5).Range.Text = "34×51="
$t.Cell(4,1).Range.Text = "88×72="
$t.Cell(4,2).Range.Text = "98×20="
$t.Cell(4,3).Range.Text = "23×10="
$t.Cell(4,4).Range.Text = "71×41="
$t.Cell(4,5).Range.Text = "97×53="
$t.Cell(5,1).Range.Text = "87×81="
$t.Cell(5,2).Range.Text = "73×96="
$t.Cell(5,3).Range.Text = "10×76="
$t.Cell(5,4).Range.Text = "17×11="
$t.Cell(5,5).Range.Text = "77×98="
$t.Cell(6,1).Range.Text = "32×76="
$t.Cell(6,2).Range.Text = "37×54="
$t.Cell(6,3).Range.Text = "41×96="
$t.Cell(6,4).Range.Text = "93×34="
$t.Cell(6,5).Range.Text = "93×62="
$t.Cell(7,1).Range.Text = "45×99="
$t.Cell(7,2).Range.Text = "80×29="
$t.Cell(7,3).Range.Text = "71×80="
$t.Cell(7,4).Range.Text = "67×43="
$t.Cell(7,5).Range.Text = "24×29="
$t.Cell(8,1).Range.Text = "82×36="
$t.Cell(8,2).Range.Text = "78×33="
$t.Cell(8,3).Range.Text = "18×75="
$t.Cell(8,4).Range.Text = "61×39="
$t.Cell(8,5).Range.Text = "98×94="
$t.Cell(9,1).Range.Text = "58×11="
$t.Cell(9,2).Range.Text = "89×39="
$t.Cell(9,3).Range.Text = "95×48="
$t.Cell(9,4).Range.Text = "18×92="
$t.Cell(9,5).Range.Text = "55×27="
$t.Cell(10,1).Range.Text = "77×41="
$t.Cell(10,2).Range.Text = "76×25="
$t.Cell(10,3).Range.Text = "62×47="
$t.Cell(10,4).Range.Text = "21×63="
$t.Cell(10,5).Range.Text = "71×25="
$t.Cell(11,1).Range.Text = "38×75="
$t.Cell(11,2).Range.Text = "21×95="
$t.Cell(11,3).Range.Text = "89×48="
$t.Cell(11,4).Range.Text = "38×13="
$t.Cell(11,5).Range.Text = "13×58="
$t.Cell(12,1).Range.Text = "39×94="
$t.Cell(12,2).Range.Text = "19×33="
$t.Cell(12,3).Range.Text = "18×67="
$t.Cell(12,4).Range.Text = "16×86="
$t.Cell(12,5).Range.Text = "38×35="
$t.Cell(13,1).Range.Text = "47×72="
$t.Cell(13,2).Range.Text = "13×58="
$t.Cell(13,3).Range.Text = "38×48="
$t.Cell(13,4).Range.Text = "24×49="
$t.Cell(13,5).Range.Text = "62×81="
$t.Cell(14,1).Range.Text = "96×47="
$t.Cell(14,2).Range.Text = "98×40="
$t.Cell(14,3).Range.Text = "72×73="
$t.Cell(14,4).Range.Text = "78×24="
$t.Cell(14,5).Range.Text = "37×21="
$t.Cell(15,1).Range.Text = "35×13="
$t.Cell(15,2).Range.Text = "47×34="
$t.Cell(15,3).Range.Text = "55×73="
$t.Cell(15,4).Range.Text = "87×36="
$t.Cell(15,5).Range.Text = "69×75="
$t.Cell(16,1).Range.Text = "27×93="
$t.Cell(16,2).Range.Text = "55×15="
$t.Cell(16,3).Range.Text = "55×67="
$t.Cell(16,4).Range.Text = "72×39="
$t.Cell(16,5).Range.Text = "15×45="
$t.Cell(17,1).Range.Text = "57×89="
$t.Cell(17,2).Range.Text = "24×65="
$t.Cell(17,3).Range.Text = "97×16="
$t.Cell(17,4).Range.Text = "74×74="
$t.Cell(17,5).Range.Text = "68×55="
$t.Cell(18,1).Range.Text = "73×22="
$t.Cell(18,2).Range.Text = "55×28="
$t.Cell(18,3).Range.Text = "86×80="
$t.Cell(18,4).Range.Text = "20×13="
$t.Cell(18,5).Range.Text = "76×94="
$t.Cell(19,1).Range.Text = "87×91="
$t.Cell(19,2).Range.Text = "99×83="
$t.Cell(19,3).Range.Text = "27×53="
$t.Cell(19,4).Range.Text = "86×58="
$t.Cell(19,5).Range.Text = "56×81="
$t.Cell(20,1).Range.Text = "42×87="
$t.Cell(20,2).Range.Text = "42×66="
$t.Cell(20,3).Range.Text = "62×48="
$t.Cell(20,4).Range.Text = "32×49="
$t.Cell(20,5).Range.Text = "44×46="
